# Insert a new weekly price record for "Berenjena" (Vega Monumental Concepción)
# as row 197, shifting the previously existing rows 197-200 down to 198-201.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 197, pushing rows 197-200 to 198-201.
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row 197 with the new data point.
$ws.Cells.Item(197, 1).Value  = 11
$ws.Cells.Item(197, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(197, 3).Value  = "Bíobío"
$ws.Cells.Item(197, 4).Value  = 45239
$ws.Cells.Item(197, 5).Value  = 8
$ws.Cells.Item(197, 6).Value  = 100112001
$ws.Cells.Item(197, 7).Value  = "Berenjena"
$ws.Cells.Item(197, 8).Value  = "Sin especificar"
$ws.Cells.Item(197, 9).Value  = "Primera"
$ws.Cells.Item(197, 10).Value = 50
$ws.Cells.Item(197, 11).Value = 12000
$ws.Cells.Item(197, 12).Value = 12000
$ws.Cells.Item(197, 13).Value = 12000
$ws.Cells.Item(197, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(197, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(197, 16).Value = 240
$ws.Cells.Item(197, 17).Value = 50
$ws.Cells.Item(197, 18).Value = "Hortaliza"
